# Fruta / hortaliza, semanal
# Weekly update: insert 4 new price rows (new week) at the top of the
# "Pera" data block (row 710), pushing the existing data down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at row 710 (existing rows 710:796 shift down to 714:800)
$ws.Range("A710:T713").Insert()

# Row 710: Packham's Triumph / Primera
$ws.Cells.Item(710,1).Value = 11
$ws.Cells.Item(710,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(710,3).Value = "Bíobío"
$ws.Cells.Item(710,4).Value = 45142
$ws.Cells.Item(710,5).Value = 8
$ws.Cells.Item(710,6).Value = "Fruta"
$ws.Cells.Item(710,7).Value = 100104
$ws.Cells.Item(710,8).Value = "Frutos de pepita"
$ws.Cells.Item(710,9).Value = 100104005
$ws.Cells.Item(710,10).Value = "Pera"
$ws.Cells.Item(710,11).Value = "Packham's Triumph"
$ws.Cells.Item(710,12).Value = "Primera"
$ws.Cells.Item(710,13).Value = 200
$ws.Cells.Item(710,14).Value = 9000
$ws.Cells.Item(710,15).Value = 10000
$ws.Cells.Item(710,16).Value = 9500
$ws.Cells.Item(710,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(710,18).Value = "Región de O'Higgins"
$ws.Cells.Item(710,19).Value = 594
$ws.Cells.Item(710,20).Value = 16

# Row 711: Packham's Triumph / Segunda
$ws.Cells.Item(711,1).Value = 11
$ws.Cells.Item(711,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(711,3).Value = "Bíobío"
$ws.Cells.Item(711,4).Value = 45142
$ws.Cells.Item(711,5).Value = 8
$ws.Cells.Item(711,6).Value = "Fruta"
$ws.Cells.Item(711,7).Value = 100104
$ws.Cells.Item(711,8).Value = "Frutos de pepita"
$ws.Cells.Item(711,9).Value = 100104005
$ws.Cells.Item(711,10).Value = "Pera"
$ws.Cells.Item(711,11).Value = "Packham's Triumph"
$ws.Cells.Item(711,12).Value = "Segunda"
$ws.Cells.Item(711,13).Value = 100
$ws.Cells.Item(711,14).Value = 8000
$ws.Cells.Item(711,15).Value = 8000
$ws.Cells.Item(711,16).Value = 8000
$ws.Cells.Item(711,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(711,18).Value = "Región de O'Higgins"
$ws.Cells.Item(711,19).Value = 500
$ws.Cells.Item(711,20).Value = 16

# Row 712: Winter Nelis / Primera
$ws.Cells.Item(712,1).Value = 11
$ws.Cells.Item(712,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(712,3).Value = "Bíobío"
$ws.Cells.Item(712,4).Value = 45142
$ws.Cells.Item(712,5).Value = 8
$ws.Cells.Item(712,6).Value = "Fruta"
$ws.Cells.Item(712,7).Value = 100104
$ws.Cells.Item(712,8).Value = "Frutos de pepita"
$ws.Cells.Item(712,9).Value = 100104005
$ws.Cells.Item(712,10).Value = "Pera"
$ws.Cells.Item(712,11).Value = "Winter Nelis"
$ws.Cells.Item(712,12).Value = "Primera"
$ws.Cells.Item(712,13).Value = 100
$ws.Cells.Item(712,14).Value = 9000
$ws.Cells.Item(712,15).Value = 10000
$ws.Cells.Item(712,16).Value = 9500
$ws.Cells.Item(712,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(712,18).Value = "Región de O'Higgins"
$ws.Cells.Item(712,19).Value = 594
$ws.Cells.Item(712,20).Value = 16

# Row 713: Winter Nelis / Segunda
$ws.Cells.Item(713,1).Value = 11
$ws.Cells.Item(713,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(713,3).Value = "Bíobío"
$ws.Cells.Item(713,4).Value = 45142
$ws.Cells.Item(713,5).Value = 8
$ws.Cells.Item(713,6).Value = "Fruta"
$ws.Cells.Item(713,7).Value = 100104
$ws.Cells.Item(713,8).Value = "Frutos de pepita"
$ws.Cells.Item(713,9).Value = 100104005
$ws.Cells.Item(713,10).Value = "Pera"
$ws.Cells.Item(713,11).Value = "Winter Nelis"
$ws.Cells.Item(713,12).Value = "Segunda"
$ws.Cells.Item(713,13).Value = 50
$ws.Cells.Item(713,14).Value = 8000
$ws.Cells.Item(713,15).Value = 8000
$ws.Cells.Item(713,16).Value = 8000
$ws.Cells.Item(713,17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(713,18).Value = "Región de O'Higgins"
$ws.Cells.Item(713,19).Value = 500
$ws.Cells.Item(713,20).Value = 16
